$d = $word.ActiveDocument

$d.Content.Find.Execute("605÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "861÷9=", 2) | Out-Null
$d.Content.Find.Execute("581÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "639÷5=", 2) | Out-Null
$d.Content.Find.Execute("911÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "871÷8=", 2) | Out-Null
$d.Content.Find.Execute("920÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "321÷9=", 2) | Out-Null
$d.Content.Find.Execute("741÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "218÷9=", 2) | Out-Null
$d.Content.Find.Execute("322÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "328÷3=", 2) | Out-Null
$d.Content.Find.Execute("558÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "749÷4=", 2) | Out-Null
$d.Content.Find.Execute("471÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "433÷5=", 2) | Out-Null
$d.Content.Find.Execute("429÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "192÷3=", 2) | Out-Null
$d.Content.Find.Execute("553÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "100÷6=", 2) | Out-Null
$d.Content.Find.Execute("391÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "145÷7=", 2) | Out-Null
$d.Content.Find.Execute("231÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "102÷2=", 2) | Out-Null
$d.Content.Find.Execute("137÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "268÷3=", 2) | Out-Null
$d.Content.Find.Execute("303÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "511÷6=", 2) | Out-Null
$d.Content.Find.Execute("248÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "810÷7=", 2) | Out-Null
$d.Content.Find.Execute("849÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "645÷8=", 2) | Out-Null
$d.Content.Find.Execute("873÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "478÷7=", 2) | Out-Null
$d.Content.Find.Execute("713÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "245÷9=", 2) | Out-Null
$d.Content.Find.Execute("226÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "166÷9=", 2) | Out-Null
$d.Content.Find.Execute("355÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "526÷6=", 2) | Out-Null
$d.Content.Find.Execute("445÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "600÷7=", 2) | Out-Null
$d.Content.Find.Execute("637÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "962÷3=", 2) | Out-Null
$d.Content.Find.Execute("294÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "527÷7=", 2) | Out-Null
$d.Content.Find.Execute("271÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "162÷8=", 2) | Out-Null
$d.Content.Find.Execute("646÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "789÷7=", 2) | Out-Null
